$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Format.WidowControl = $false
}
